{"js": "// \"Transfert de messages et de fichiers\" paragraph: replace the placeholder\n// ellipsis (\"\u2026\u2026\u2026\u2026\u2026\u2026\u2026\u2026..\") that stands in for the maximum file size with the\n// actual value \"4 Go\", leaving the rest of the sentence untouched.\nconst body = context.document.body;\n\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst placeholder = \"\u2026\u2026\u2026\u2026\u2026\u2026\u2026\u2026..\";\nlet paragraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(placeholder) !== -1) {\n    paragraph = paragraphs.items[i];\n    break;\n  }\n}\nif (!paragraph) {\n  throw new Error(\"Could not find the paragraph containing the placeholder ellipsis.\");\n}\n\n// Scope the search to this paragraph so we only ever touch this one spot,\n// even though the placeholder text is already unique document-wide.\nconst results = paragraph.search(placeholder, { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Placeholder ellipsis text not found in the target paragraph.\");\n}\n\nresults.items[0].insertText(\"4 Go\", \"Replace\");\nawait context.sync();\n", "ps1": "# \"Transfert de messages et de fichiers\" paragraph: replace the placeholder\n# ellipsis (\"\u2026\u2026\u2026\u2026\u2026\u2026\u2026\u2026..\") that stands in for the maximum file size with the\n# actual value \"4 Go\", leaving the rest of the sentence untouched.\n\n$d = $word.ActiveDocument\n\n# Locate the paragraph that talks about the maximum file size \u2014 this is more\n# robust than assuming a fixed paragraph index.\n$target = $null\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text -like \"*taille maximale des fichiers*\") {\n        $target = $p.Range\n        break\n    }\n}\nif ($target -eq $null) {\n    throw \"Could not find the paragraph containing 'taille maximale des fichiers'.\"\n}\n\n# Scope the Find/Replace to that paragraph's range so we only ever touch this\n# one spot, even though the placeholder text is already unique document-wide.\n$find = $target.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"\u2026\u2026\u2026\u2026\u2026\u2026\u2026\u2026..\"\n$find.Replacement.Text = \"4 Go\"\n\n$replaced = $find.Execute(\n    $find.Text,        # FindText\n    $false,            # MatchCase\n    $false,            # MatchWholeWord\n    $false,            # MatchWildcards\n    $false,            # MatchSoundsLike\n    $false,            # MatchAllWordForms\n    $true,             # Forward\n    1,                 # Wrap (wdFindContinue)\n    $false,            # Format\n    $find.Replacement.Text,  # ReplaceWith\n    2                  # Replace (wdReplaceAll)\n)\n\nif (-not $replaced) {\n    throw \"Placeholder ellipsis text not found in the target paragraph.\"\n}\n"}
